# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This writes the recalculated "K" column (column G) values for rows 2-56.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 1
    4  = 2
    5  = 1
    6  = 2
    7  = 0
    8  = 1
    9  = 0
    10 = 1
    11 = 1
    12 = 2
    13 = 2
    14 = 0
    15 = 1
    16 = 1
    17 = 2
    18 = 1
    19 = 0
    20 = 1
    21 = 1
    22 = 2
    23 = 1
    24 = 0
    25 = 1
    26 = 1
    27 = 1
    28 = 2
    29 = 0
    30 = 0
    31 = 2
    32 = 2
    33 = 0
    34 = 1
    35 = 3
    36 = 1
    37 = 2
    38 = 2
    39 = 2
    40 = 0
    41 = 1
    42 = 1
    43 = 4
    44 = 2
    45 = 3
    46 = 2
    47 = 3
    48 = 1
    49 = 3
    50 = 1
    51 = 5
    52 = 1
    53 = 4
    54 = 1
    55 = 4
    56 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}

Write-Host "Updated column G (K) for rows 2-56"
